# Update the cached "today" date and the slide-number placeholder glyph
# on the slide master and on every slide layout.
#
#   datetimeFigureOut field text: "12/12/2020" -> "29/12/2020"
#   slidenum field text:          "<Nº>"       -> "<#>"
#
# ppPlaceholderDate = 16, ppPlaceholderSlideNumber = 13 (MsoPlaceholderType)

$p = $ppt.ActivePresentation

$newDate = "29/12/2020"
$newSlideNum = [string]([char]0x2039) + [string]([char]0x23) + [string]([char]0x203A)

function Update-HeaderFooterShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)

        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch {}

        if ($phType -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        } elseif ($phType -eq 13) {
            $shp.TextFrame.TextRange.Text = $newSlideNum
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-HeaderFooterShapes $master

# Every slide layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-HeaderFooterShapes $layout
}
